$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
Write-Output $excel.ActiveWindow.Zoom()
